# Add a new "2022" column (column R) to the table, mirroring the existing
# 2021 column (Q), then leave the active selection on S4 (one cell to the
# right of the newly added data), matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the whole 2021 column (values + formatting) into the new column so
# the new cells pick up the same number formatting / borders / fonts as
# their neighbours in column Q.
$ws.Range("Q4:Q6").Copy($ws.Range("R4:R6"))

# Overwrite the copied values with the new 2022 figures.
$ws.Range("R4").Value = 2022
$ws.Range("R5").Value = 8.6821914120339212
$ws.Range("R6").Value = 12.221423436376707

# Match the saved selection/active cell state recorded in the workbook.
$ws.Range("S4").Select()
